$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.001.43'
$ws.Range('E2').Value = '  -1.03%  '

$ws.Range('D3').Value = '3.420.82'
$ws.Range('E3').Value = '  +3.38%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '''255.85'
$ws.Range('E5').Value = '  -0.36%  '

$ws.Range('D6').Value = '''656.36'
$ws.Range('E6').Value = '  +4.58%  '

$ws.Range('E7').Value = '  +2.48%  '

$ws.Range('D8').Value = '''0.439'
$ws.Range('E8').Value = '  +6.08%  '

$ws.Range('D9').Value = '''1.09'
$ws.Range('E9').Value = '  +7.41%  '

$ws.Range('D11').Value = '3.420.13'
$ws.Range('E11').Value = '  +3.41%  '

$ws.Range('D12').Value = '''0.213'
$ws.Range('E12').Value = '  +4.60%  '

$ws.Range('D13').Value = '''42.15'
$ws.Range('E13').Value = '  +1.77%  '

$ws.Range('D14').Value = '''6.44'
$ws.Range('E14').Value = '  +19.25%  '

$ws.Range('D15').Value = '''0.0000261'
$ws.Range('E15').Value = '  +3.30%  '

$ws.Range('D16').Value = '97.708.97'
$ws.Range('E16').Value = '  -1.10%  '

$ws.Range('D17').Value = '4.066.61'
$ws.Range('E17').Value = '  +3.39%  '

$ws.Range('D18').Value = '''8.68'
$ws.Range('E18').Value = '  +34.68%  '

$ws.Range('D19').Value = '3.427.55'
$ws.Range('E19').Value = '  +3.73%  '

$ws.Range('D20').Value = '''17.73'
$ws.Range('E20').Value = '  +13.07%  '

$ws.Range('D21').Value = '''0.493'
$ws.Range('E21').Value = '  +46.66%  '

$ws.Range('D22').Value = '''3.48'
$ws.Range('E22').Value = '  +0.05%  '

$ws.Range('D23').Value = '''10.77'
$ws.Range('E23').Value = '  +13.69%  '

$ws.Range('D24').Value = '''513.03'
$ws.Range('E24').Value = '  +5.09%  '

$ws.Range('D25').Value = '''0.0000208'
$ws.Range('E25').Value = '  +2.04%  '

$ws.Range('D26').Value = '''6.24'
$ws.Range('E26').Value = '  +7.89%  '

$ws.Range('E27').Value = '  +11.27%  '

$ws.Range('D28').Value = '''12.97'
$ws.Range('E28').Value = '  +5.81%  '

$ws.Range('D29').Value = '3.605.02'
$ws.Range('E29').Value = '  +3.39%  '

$ws.Range('E30').Value = '  +0.73%  '

$ws.Range('D31').Value = '''0.202'
$ws.Range('E31').Value = '  +5.18%  '

$ws.Range('D32').Value = '''11.48'
$ws.Range('E32').Value = '  +8.23%  '

$ws.Range('D33').Value = '''1.00'
$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('D34').Value = '''1.00'
$ws.Range('E34').Value = '  -0.01%  '

$ws.Range('D35').Value = '''0.575'
$ws.Range('E35').Value = '  +19.05%  '

$ws.Range('E36').Value = '  +6.39%  '

$ws.Range('D37').Value = '''2.30'
$ws.Range('E37').Value = '  +16.68%  '

$ws.Range('D38').Value = '''7.87'
$ws.Range('E38').Value = '  +6.60%  '

$ws.Range('E39').Value = '  +14.45%  '

$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '''529.93'
$ws.Range('E40').Value = '  +6.34%  '

$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.155'
$ws.Range('E41').Value = '  +1.44%  '

$ws.Range('E42').Value = '  +0.10%  '

$ws.Range('D43').Value = '''0.870'
$ws.Range('E43').Value = '  +10.26%  '

$ws.Range('D44').Value = '''3.70'
$ws.Range('E44').Value = '  -3.82%  '

$ws.Range('D45').Value = '''0.0420'
$ws.Range('E45').Value = '  +21.99%  '

$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').Value = '''3.34'
$ws.Range('E46').Value = '  +3.37%  '

$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '''5.54'
$ws.Range('E47').Value = '  +15.49%  '

$ws.Range('D48').Value = '''8.34'
$ws.Range('E48').Value = '  +13.21%  '

$ws.Range('B49').Value = 'ImmutableX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D49').Value = '''1.60'
$ws.Range('E49').Value = '  +13.76%  '

$ws.Range('B50').Value = 'USDe'
$ws.Range('C50').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D50').Value = '''1.00'
$ws.Range('E50').Value = '  +0.02%  '

$ws.Range('D51').Value = '''2.08'
$ws.Range('E51').Value = '  +5.84%  '
